# Applies the "Added a few more slots" edit to blazing-mammoth:
#  1. Insert a new "Meta description: ..." paragraph right after the H1 title.
#  2. Remove the duplicated bold title paragraph near the end of the document.
#  3. Replace the text of the final (italic) paragraph with the image-prompt text.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Step 1: insert the "Meta description" paragraph after paragraph 1 ---------
$titlePara = $d.Paragraphs.Item(1)
[void]$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaXml = '<w:p ' + $wNs + '>' +
    '<w:r/>' +
    '<w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>' +
    '<w:r><w:t>: Read our unbiased review of Blazing Mammoth, a 5-reel virtual slot machine with prehistoric theme. Play it for free and enjoy unique game features!</w:t></w:r>' +
    '</w:p>'
[void]$metaPara.Range.InsertXML($metaXml)

# --- Step 2: delete the duplicated bold "Play Blazing Mammoth Free ..." -------
#     paragraph that used to sit right before the italic meta-description       
#     paragraph near the end of the document.                                   
$found = $false
$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($i -gt 1 -and $p.Range.Text.StartsWith("Play Blazing Mammoth Free | A Prehistoric Themed Slot Game")) {
        $p.Range.Delete()
        $found = $true
        break
    }
}

# --- Step 3: replace the text of the last paragraph with the image prompt -----
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$promptText = 'Prompt: Create a feature image fitting the game "Blazing Mammoth". The image should be in cartoon style featuring a happy Maya warrior with glasses. Specifically, the image should have the followings: - The Maya warrior should be standing in front of a big Blazing Mammoth with a smile on his face, representing the fun and exciting aspect of the game. - The cartoon style should be colorful and vibrant, making it appealing and eye-catching to potential players. - The glasses on the Maya warrior should be distinctive and represent that the game is modern and easy to access. - It should be designed in a way that it represents the prehistoric era and the concept of the game.'
$promptXml = '<w:p ' + $wNs + '>' +
    '<w:r/>' +
    '<w:r><w:rPr><w:i/></w:rPr><w:t>' + $promptText + '</w:t></w:r>' +
    '</w:p>'
[void]$lastPara.Range.InsertXML($promptXml)

Write-Host "done; found-and-removed duplicate title paragraph: $found"
